{"js": "// Update the three-digit-by-one-digit multiplication prompts throughout\n// the worksheet table. Each \"old\u00d7old=\" string is unique in the document,\n// so a direct search + replace per pair is unambiguous.\nconst replacements = [\n  [\"121\u00d72=\", \"465\u00d74=\"],\n  [\"829\u00d76=\", \"405\u00d75=\"],\n  [\"703\u00d75=\", \"151\u00d77=\"],\n  [\"985\u00d78=\", \"541\u00d73=\"],\n  [\"952\u00d72=\", \"446\u00d75=\"],\n  [\"392\u00d79=\", \"102\u00d74=\"],\n  [\"416\u00d72=\", \"519\u00d73=\"],\n  [\"438\u00d79=\", \"415\u00d76=\"],\n  [\"708\u00d78=\", \"304\u00d77=\"],\n  [\"341\u00d73=\", \"403\u00d75=\"],\n  [\"323\u00d79=\", \"679\u00d79=\"],\n  [\"511\u00d74=\", \"612\u00d76=\"],\n  [\"904\u00d77=\", \"106\u00d76=\"],\n  [\"998\u00d75=\", \"565\u00d75=\"],\n  [\"947\u00d75=\", \"120\u00d72=\"],\n  [\"865\u00d79=\", \"359\u00d76=\"],\n  [\"742\u00d79=\", \"588\u00d74=\"],\n  [\"858\u00d75=\", \"909\u00d79=\"],\n  [\"639\u00d74=\", \"794\u00d76=\"],\n  [\"877\u00d79=\", \"531\u00d77=\"],\n  [\"201\u00d72=\", \"830\u00d76=\"],\n  [\"366\u00d73=\", \"150\u00d78=\"],\n  [\"491\u00d77=\", \"838\u00d78=\"],\n  [\"883\u00d74=\", \"866\u00d75=\"],\n  [\"605\u00d77=\", \"803\u00d72=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the three-digit-by-one-digit multiplication prompts throughout\n# the worksheet table. Each \"old\u00d7old=\" string is unique in the document,\n# so a direct Find/Replace per pair is unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"121\u00d72=\", \"465\u00d74=\"),\n    @(\"829\u00d76=\", \"405\u00d75=\"),\n    @(\"703\u00d75=\", \"151\u00d77=\"),\n    @(\"985\u00d78=\", \"541\u00d73=\"),\n    @(\"952\u00d72=\", \"446\u00d75=\"),\n    @(\"392\u00d79=\", \"102\u00d74=\"),\n    @(\"416\u00d72=\", \"519\u00d73=\"),\n    @(\"438\u00d79=\", \"415\u00d76=\"),\n    @(\"708\u00d78=\", \"304\u00d77=\"),\n    @(\"341\u00d73=\", \"403\u00d75=\"),\n    @(\"323\u00d79=\", \"679\u00d79=\"),\n    @(\"511\u00d74=\", \"612\u00d76=\"),\n    @(\"904\u00d77=\", \"106\u00d76=\"),\n    @(\"998\u00d75=\", \"565\u00d75=\"),\n    @(\"947\u00d75=\", \"120\u00d72=\"),\n    @(\"865\u00d79=\", \"359\u00d76=\"),\n    @(\"742\u00d79=\", \"588\u00d74=\"),\n    @(\"858\u00d75=\", \"909\u00d79=\"),\n    @(\"639\u00d74=\", \"794\u00d76=\"),\n    @(\"877\u00d79=\", \"531\u00d77=\"),\n    @(\"201\u00d72=\", \"830\u00d76=\"),\n    @(\"366\u00d73=\", \"150\u00d78=\"),\n    @(\"491\u00d77=\", \"838\u00d78=\"),\n    @(\"883\u00d74=\", \"866\u00d75=\"),\n    @(\"605\u00d77=\", \"803\u00d72=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $r = $d.Content\n    $r.Find.ClearFormatting()\n    $r.Find.Replacement.ClearFormatting()\n    $r.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
